$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bottle Results")
$ws.Range("A11:AB13").Font.ThemeColor = 1
Write-Host "done"
